# Rename the inline picture shapes (the Pearson/BTec logos living in the
# document's headers/footers) so each one's wp:docPr/pic:cNvPr "name"
# attribute matches the target file name:
#   - footer (default / primary)    : PearsonLogo  image2.png -> image1.png
#   - footer (first page)           : PearsonLogo  image2.png -> image1.png
#   - header (first page)           : BTec logo    image1.jpg -> image2.jpg
#
# NOTE: InlineShape objects returned directly off a HeaderFooter's
# .Range.InlineShapes collection can resolve to a stale/invalid anchor in
# this host when the footer/header story hasn't been touched yet, which
# makes a direct ".Name = ..." assignment throw. Selecting the shape first
# and then re-fetching it from $word.Selection.InlineShapes works reliably,
# so every rename below goes through that route.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-InlineLogo($headerFooter, $newName) {
    $shape = $headerFooter.Range.InlineShapes.Item(1)
    $shape.Select()
    $word.Selection.InlineShapes.Item(1).Name = $newName
}

# Footer, default/primary (word/footer2.xml) - PearsonLogo, id="2"
Rename-InlineLogo $sec.Footers.Item(1) "image1.png"

# Footer, first page (word/footer1.xml) - PearsonLogo, id="3"
Rename-InlineLogo $sec.Footers.Item(2) "image1.png"

# Header, first page (word/header1.xml) - BTec_Logo-Orange, id="1"
Rename-InlineLogo $sec.Headers.Item(2) "image2.jpg"
